$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "1.002"); the source
# data stores these as text, so force a Text number format before assigning,
# otherwise Excel auto-converts the literal to a Number cell.
$textCells = @(
    "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13",
    "D14", "D15", "D16", "D19", "D20", "D21", "D22", "D23", "D25", "D26",
    "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38",
    "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48",
    "D49", "D50", "D51"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# New cell values (row-major, matching the updated cryptos snapshot).
$updates = [ordered]@{
    "D2" = "23.536.22"
    "E2" = "  +0.72%  "
    "D3" = "1.652.32"
    "E3" = "  +1.31%  "
    "D4" = "1.002"
    "E4" = "  +0.31%  "
    "D5" = "1.002"
    "E5" = "  +0.42%  "
    "D6" = "300.06"
    "E6" = "  -1.01%  "
    "D7" = "0.3800"
    "E7" = "  +1.26%  "
    "D8" = "50.82"
    "E8" = "  -1.66%  "
    "D9" = "0.3561"
    "E9" = "  -0.22%  "
    "D10" = "0.08112"
    "E10" = "  -0.79%  "
    "D11" = "1.226"
    "E11" = "  -0.45%  "
    "D12" = "1.002"
    "E12" = "  +0.32%  "
    "D13" = "22.07"
    "E13" = "  -0.88%  "
    "D14" = "6.412"
    "E14" = "  -1.16%  "
    "D15" = "7.403"
    "E15" = "  +1.14%  "
    "D16" = "0.00001203"
    "E16" = "  -1.52%  "
    "D17" = "1.656.17"
    "E17" = "  +2.25%  "
    "E18" = "  +1.77%  "
    "D19" = "0.06992"
    "E19" = "  +0.78%  "
    "D20" = "6.791"
    "E20" = "  +1.37%  "
    "D21" = "17.52"
    "E21" = "  +0.54%  "
    "D22" = "1.001"
    "E22" = "  +0.33%  "
    "D23" = "12.66"
    "E23" = "  +1.78%  "
    "D24" = "23.562.50"
    "E24" = "  +0.90%  "
    "D25" = "2.491"
    "E25" = "  -1.11%  "
    "D26" = "2.953"
    "E26" = "  -4.39%  "
    "E27" = "  -0.29%  "
    "D28" = "151.88"
    "E28" = "  -0.88%  "
    "D29" = "5.236"
    "E29" = "  +1.12%  "
    "D30" = "133.39"
    "E30" = "  -0.31%  "
    "D31" = "1.838.33"
    "E31" = "  +2.08%  "
    "D32" = "6.989"
    "E32" = "  +5.59%  "
    "D33" = "2.161"
    "E33" = "  +6.84%  "
    "D34" = "11.77"
    "E34" = "  +1.50%  "
    "D35" = "1.039"
    "E35" = "  -5.11%  "
    "D36" = "0.02744"
    "E36" = "  +0.13%  "
    "D37" = "0.08728"
    "E37" = "  -0.32%  "
    "B38" = "InternetComputer(DFINITY)"
    "C38" = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
    "D38" = "5.990"
    "E38" = "  +0.91%  "
    "B39" = "Algorand"
    "C39" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
    "D39" = "0.2450"
    "E39" = "  -1.28%  "
    "D40" = "13.20"
    "E40" = "  +5.23%  "
    "D41" = "0.06870"
    "E41" = "  -1.13%  "
    "D42" = "0.6943"
    "E42" = "  -0.04%  "
    "D43" = "1.321"
    "E43" = "  -0.49%  "
    "D44" = "15.71"
    "E44" = "  +0.86%  "
    "D45" = "0.6450"
    "E45" = "  +0.85%  "
    "D46" = "1.001"
    "D47" = "2.268"
    "E47" = "  -1.01%  "
    "D48" = "3.927"
    "E48" = "  -0.77%  "
    "D49" = "0.07876"
    "E49" = "  -0.55%  "
    "D50" = "126.77"
    "E50" = "  -0.64%  "
    "D51" = "1.178"
    "E51" = "  +0.17%  "
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
